$p = $ppt.ActivePresentation

function ComRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# The deck's design ("Integral" / Red Violet) is swapped for the default
# "Office Theme" colour palette - i.e. applying a new Office theme to the
# slide master (ppt/theme/theme1.xml), which drives every slide's colours.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = ComRGB 0x00 0x00 0x00   # dk1      -> 000000
$tcs.Item(2).RGB  = ComRGB 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = ComRGB 0x44 0x54 0x6A   # dk2      -> 44546A
$tcs.Item(4).RGB  = ComRGB 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = ComRGB 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = ComRGB 0xED 0x7D 0x31   # accent2  -> ED7D31
$tcs.Item(7).RGB  = ComRGB 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = ComRGB 0xFF 0xC0 0x00   # accent4  -> FFC000
$tcs.Item(9).RGB  = ComRGB 0x44 0x72 0xC4   # accent5  -> 4472C4
$tcs.Item(10).RGB = ComRGB 0x70 0xAD 0x47   # accent6  -> 70AD47
$tcs.Item(11).RGB = ComRGB 0x05 0x63 0xC1   # hlink    -> 0563C1
$tcs.Item(12).RGB = ComRGB 0x95 0x4F 0x72   # folHlink -> 954F72
